$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 374-375 (shifts existing data, incl. rows 374-490, down by 2
# to 376-492) to make room for the newest week of price data at the top of the table.
$ws.Range("A374:A375").EntireRow.Insert()

# Row 374: Apio, Americana (o), Primera
$ws.Cells.Item(374, 1).Value = 11
$ws.Cells.Item(374, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(374, 3).Value = "Bíobío"
$ws.Cells.Item(374, 4).Value = 45093
$ws.Cells.Item(374, 5).Value = 8
$ws.Cells.Item(374, 6).Value = 100112017
$ws.Cells.Item(374, 7).Value = "Apio"
$ws.Cells.Item(374, 8).Value = "Americana (o)"
$ws.Cells.Item(374, 9).Value = "Primera"
$ws.Cells.Item(374, 10).Value = 100
$ws.Cells.Item(374, 11).Value = 7500
$ws.Cells.Item(374, 12).Value = 8000
$ws.Cells.Item(374, 13).Value = 7750
$ws.Cells.Item(374, 14).Value = "$/docena de matas"
$ws.Cells.Item(374, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(374, 16).Value = 1292
$ws.Cells.Item(374, 17).Value = 6
$ws.Cells.Item(374, 18).Value = "Hortaliza"

# Row 375: Apio, Americana (o), Segunda
$ws.Cells.Item(375, 1).Value = 11
$ws.Cells.Item(375, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(375, 3).Value = "Bíobío"
$ws.Cells.Item(375, 4).Value = 45093
$ws.Cells.Item(375, 5).Value = 8
$ws.Cells.Item(375, 6).Value = 100112017
$ws.Cells.Item(375, 7).Value = "Apio"
$ws.Cells.Item(375, 8).Value = "Americana (o)"
$ws.Cells.Item(375, 9).Value = "Segunda"
$ws.Cells.Item(375, 10).Value = 50
$ws.Cells.Item(375, 11).Value = 6500
$ws.Cells.Item(375, 12).Value = 6500
$ws.Cells.Item(375, 13).Value = 6500
$ws.Cells.Item(375, 14).Value = "$/docena de matas"
$ws.Cells.Item(375, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(375, 16).Value = 1083
$ws.Cells.Item(375, 17).Value = 6
$ws.Cells.Item(375, 18).Value = "Hortaliza"
